$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Wnt11"
$ws.Range("C2").Value = "Fzd7"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.062425
$ws.Range("H2").Value = 0.187275
$ws.Range("I2").Value = 0.00296197839880675
$ws.Range("J2").Value = 0.00296197839880675
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 2.418393
$ws.Range("N2").Value = 7.255179
$ws.Range("O2").Value = 0.0919828589765645
$ws.Range("P2").Value = 0.0919828589765645
$ws.Range("Q2").Value = 0.150968183025
$ws.Range("R2").Value = 1.358713647225
$ws.Range("S2").Value = 0.0002724512413490716
$ws.Range("T2").Value = 0.0002724512413490716

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Wnt11"
$ws.Range("C3").Value = "Fzd7"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.062425
$ws.Range("H3").Value = 0.187275
$ws.Range("I3").Value = 0.00296197839880675
$ws.Range("J3").Value = 0.00296197839880675
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 10.11799233333333
$ws.Range("N3").Value = 30.353977
$ws.Range("O3").Value = 0.3848348311969811
$ws.Range("P3").Value = 0.3848348311969811
$ws.Range("Q3").Value = 0.6316156714083333
$ws.Range("R3").Value = 5.684541042675
$ws.Range("S3").Value = 0.0011398724571139
$ws.Range("T3").Value = 0.0011398724571139

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Wnt11"
$ws.Range("C4").Value = "Fzd7"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.062425
$ws.Range("H4").Value = 0.187275
$ws.Range("I4").Value = 0.00296197839880675
$ws.Range("J4").Value = 0.00296197839880675
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 13.75539366666667
$ws.Range("N4").Value = 41.266181
$ws.Range("O4").Value = 0.5231823098264544
$ws.Range("P4").Value = 0.5231823098264544
$ws.Range("Q4").Value = 0.8586804496416668
$ws.Range("R4").Value = 7.728124046775
$ws.Range("S4").Value = 0.001549654700343778
$ws.Range("T4").Value = 0.001549654700343779

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Wnt11"
$ws.Range("C5").Value = "Fzd7"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 20.223983
$ws.Range("H5").Value = 60.671949
$ws.Range("I5").Value = 0.9595995319797346
$ws.Range("J5").Value = 0.9595995319797347
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 2.418393
$ws.Range("N5").Value = 7.255179
$ws.Range("O5").Value = 0.0919828589765645
$ws.Range("P5").Value = 0.0919828589765645
$ws.Range("Q5").Value = 48.909538919319
$ws.Range("R5").Value = 440.185850273871
$ws.Range("S5").Value = 0.08826670842406922
$ws.Range("T5").Value = 0.08826670842406922

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Wnt11"
$ws.Range("C6").Value = "Fzd7"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 20.223983
$ws.Range("H6").Value = 60.671949
$ws.Range("I6").Value = 0.9595995319797346
$ws.Range("J6").Value = 0.9595995319797347
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 10.11799233333333
$ws.Range("N6").Value = 30.353977
$ws.Range("O6").Value = 0.3848348311969811
$ws.Range("P6").Value = 0.3848348311969811
$ws.Range("Q6").Value = 204.6261049434637
$ws.Range("R6").Value = 1841.634944491173
$ws.Range("S6").Value = 0.3692873239061232
$ws.Range("T6").Value = 0.3692873239061233

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Wnt11"
$ws.Range("C7").Value = "Fzd7"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 20.223983
$ws.Range("H7").Value = 60.671949
$ws.Range("I7").Value = 0.9595995319797346
$ws.Range("J7").Value = 0.9595995319797347
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 13.75539366666667
$ws.Range("N7").Value = 41.266181
$ws.Range("O7").Value = 0.5231823098264544
$ws.Range("P7").Value = 0.5231823098264544
$ws.Range("Q7").Value = 278.1888476729744
$ws.Range("R7").Value = 2503.699629056769
$ws.Range("S7").Value = 0.5020454996495421
$ws.Range("T7").Value = 0.5020454996495423

# Row 8
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Wnt11"
$ws.Range("C8").Value = "Fzd7"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.7890326666666666
$ws.Range("H8").Value = 2.367098
$ws.Range("I8").Value = 0.0374384896214586
$ws.Range("J8").Value = 0.03743848962145861
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 2.418393
$ws.Range("N8").Value = 7.255179
$ws.Range("O8").Value = 0.0919828589765645
$ws.Range("P8").Value = 0.0919828589765645
$ws.Range("Q8").Value = 1.908191077838
$ws.Range("R8").Value = 17.173719700542
$ws.Range("S8").Value = 0.0034436993111462
$ws.Range("T8").Value = 0.003443699311146201

# Row 9
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Wnt11"
$ws.Range("C9").Value = "Fzd7"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.7890326666666666
$ws.Range("H9").Value = 2.367098
$ws.Range("I9").Value = 0.0374384896214586
$ws.Range("J9").Value = 0.03743848962145861
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 10.11799233333333
$ws.Range("N9").Value = 30.353977
$ws.Range("O9").Value = 0.3848348311969811
$ws.Range("P9").Value = 0.3848348311969811
$ws.Range("Q9").Value = 7.983426472082888
$ws.Range("R9").Value = 71.850838248746
$ws.Range("S9").Value = 0.01440763483374395
$ws.Range("T9").Value = 0.01440763483374395

# Row 10
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Wnt11"
$ws.Range("C10").Value = "Fzd7"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.7890326666666666
$ws.Range("H10").Value = 2.367098
$ws.Range("I10").Value = 0.0374384896214586
$ws.Range("J10").Value = 0.03743848962145861
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 13.75539366666667
$ws.Range("N10").Value = 41.266181
$ws.Range("O10").Value = 0.5231823098264544
$ws.Range("P10").Value = 0.5231823098264544
$ws.Range("Q10").Value = 10.85345494585978
$ws.Range("R10").Value = 97.681094512738
$ws.Range("S10").Value = 0.01958715547656845
$ws.Range("T10").Value = 0.01958715547656845
